$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/17/2025  Through  2/23/2025"

# --- Numeric cell updates in the crime statistics table (rows 16-28) ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 15
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = 36.363636363636
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -42.307692307692
$ws.Range("N16").Value = -87.903225806451

$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 24
$ws.Range("J17").Value = 12
$ws.Range("K17").Value = 100
$ws.Range("L17").Value = 140
$ws.Range("M17").Value = 166.666666666667
$ws.Range("N17").Value = 71.428571428571

$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 116.666666666667
$ws.Range("I18").Value = 33
$ws.Range("J18").Value = 14
$ws.Range("K18").Value = 135.714285714286
$ws.Range("L18").Value = 17.857142857142
$ws.Range("M18").Value = 73.684210526315
$ws.Range("N18").Value = -85.897435897435

$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -25.714285714285
$ws.Range("I19").Value = 51
$ws.Range("J19").Value = 55
$ws.Range("K19").Value = -7.272727272727
$ws.Range("L19").Value = -25
$ws.Range("M19").Value = -17.741935483871
$ws.Range("N19").Value = -63.043478260869

$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 12.5
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 14
$ws.Range("L20").Value = -4.545454545454
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -96.391752577319

$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -12.5
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 61
$ws.Range("H21").Value = 8.196721311475
$ws.Range("I21").Value = 146
$ws.Range("J21").Value = 107
$ws.Range("K21").Value = 36.448598130841
$ws.Range("L21").Value = 5.797101449275
$ws.Range("M21").Value = 6.569343065693
$ws.Range("N21").Value = -86.654478976234

$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 25
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = 40
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = 75

$ws.Range("C24").Value = 51
$ws.Range("E24").Value = 34.210526315789
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 139
$ws.Range("H24").Value = 34.532374100719
$ws.Range("I24").Value = 324
$ws.Range("J24").Value = 246
$ws.Range("K24").Value = 31.707317073170
$ws.Range("L24").Value = 31.707317073170
$ws.Range("M24").Value = 205.660377358491

$ws.Range("C25").Value = 41
$ws.Range("D25").Value = 29
$ws.Range("E25").Value = 41.379310344827
$ws.Range("F25").Value = 154
$ws.Range("G25").Value = 113
$ws.Range("H25").Value = 36.283185840708
$ws.Range("I25").Value = 267
$ws.Range("J25").Value = 192
$ws.Range("K25").Value = 39.0625
$ws.Range("L25").Value = 50

$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -77.777777777777
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 94.117647058823
$ws.Range("I26").Value = 49
$ws.Range("J26").Value = 32
$ws.Range("K26").Value = 53.125
$ws.Range("L26").Value = 63.333333333333
$ws.Range("M26").Value = 40

$ws.Range("F27").Value = 1
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 100

$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 200

# --- Cells converted from numeric counts to the text placeholders "0" / "***.*" ---
# Value is set with a leading apostrophe to force text, then the number format/style
# is copied from an existing placeholder cell in the same row so the look matches
# the other "N/A" style cells in the table.
$ws.Range("F15").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("G15").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H15").Value = "'***.*"
$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial(-4122)

$ws.Range("G27").Value = "'0"
$ws.Range("D27").Copy()
$ws.Range("G27").PasteSpecial(-4122)

$ws.Range("H27").Value = "'***.*"
$ws.Range("E27").Copy()
$ws.Range("H27").PasteSpecial(-4122)

$ws.Range("C31").Value = "'0"
$ws.Range("D31").Copy()
$ws.Range("C31").PasteSpecial(-4122)

$excel.CutCopyMode = 0
